$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Optical_Power")

# Update existing "Pendiente ADM" placeholders in column E with the real OT codes
$ws.Range("E9").Value = "ICD30399137"
$ws.Range("E10").Value = "ICD30399168"
$ws.Range("E11").Value = "ICD30399224"
$ws.Range("E13").Value = "ICD30399267"
$ws.Range("E14").Value = "ICD30399377"
$ws.Range("E15").Value = "ICD30388090"
$ws.Range("E16").Value = "ICD30398505"

# Append a new row of data (row 17)
# Leading apostrophes force Excel to store these as plain text (matching the
# source data) instead of auto-converting to a number / date.
$ws.Range("A17").Value = "'6946"
$ws.Range("B17").Value = "'8/14/2025"
$ws.Range("C17").Value = "SAENZ AV. 1161"
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = "ICD30398322 "
$ws.Range("F17").Value = "Optical Power"
$ws.Range("G17").Value = "Pendiente"
$ws.Range("H17").Value = "Tendido a baja altura"
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = '{"direccionesNormalizadas": [{"altura": 1161, "cod_calle": 20002, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.416595", "y": "-34.653750"}, "direccion": "SAENZ AV. 1161, CABA", "nombre_calle": "SAENZ AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}, {"altura": 1161, "cod_calle": 20004, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.382941", "y": "-34.602931"}, "direccion": "SAENZ PEÑA, ROQUE, PRES. DIAGONAL NORTE AV. 1161, CABA", "nombre_calle": "SAENZ PEÑA, ROQUE, PRES. DIAGONAL NORTE AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K17").Value = -58.416595
$ws.Range("L17").Value = -34.65375
$ws.Range("M17").Value = "San Telmo"
$ws.Range("N17").Value = "Capital Sur"
